$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ D = "36.793.90"; E = "  -0.92%  " }
    3 = @{ D = "2.079.61"; E = "  +1.44%  " }
    4 = @{ D = $null; E = "  -0.24%  " }
    5 = @{ D = "244.89"; E = "  -1.24%  " }
    6 = @{ D = $null; E = "  -1.93%  " }
    7 = @{ D = $null; E = "  -0.02%  " }
    8 = @{ D = "54.35"; E = "  -6.32%  " }
    9 = @{ D = "59.66"; E = "  -0.36%  " }
    10 = @{ D = $null; E = "  -4.04%  " }
    11 = @{ D = "0.0761"; E = "  -1.78%  " }
    12 = @{ D = $null; E = "  +1.33%  " }
    13 = @{ D = "14.97"; E = "  -5.76%  " }
    14 = @{ D = "0.883"; E = "  +2.76%  " }
    15 = @{ D = "2.382.90"; E = "  +1.45%  " }
    16 = @{ D = "5.49"; E = "  -3.53%  " }
    17 = @{ D = "2.050.91"; E = "  -0.03%  " }
    18 = @{ D = "36.713.22"; E = "  -1.00%  " }
    19 = @{ D = $null; E = "  -3.77%  " }
    20 = @{ D = "72.71"; E = "  -2.83%  " }
    21 = @{ D = "0.0₃0879"; E = "  -1.13%  " }
    22 = @{ D = "5.44"; E = "  +1.85%  " }
    23 = @{ D = "237.81"; E = "  +0.32%  " }
    24 = @{ D = "1.00"; E = "  +0.14%  " }
    25 = @{ D = "2.41"; E = "  -1.88%  " }
    26 = @{ D = "9.82"; E = "  +3.64%  " }
    27 = @{ D = $null; E = "  -0.36%  " }
    28 = @{ D = "167.32"; E = "  -1.37%  " }
    29 = @{ D = "20.65"; E = "  +3.15%  " }
    30 = @{ D = $null; E = "  -1.09%  " }
    31 = @{ D = "5.27"; E = "  +10.12%  " }
    32 = @{ D = "1.19"; E = "  +6.77%  " }
    33 = @{ D = $null; E = "  +4.25%  " }
    34 = @{ D = $null; E = "  -1.44%  " }
    35 = @{ D = "2.37"; E = "  +5.46%  " }
    36 = @{ D = $null; E = "  +0.09%  " }
    37 = @{ D = $null; E = "  +4.02%  " }
    38 = @{ D = "0.0835"; E = "  -6.74%  " }
    39 = @{ D = $null; E = "  -4.60%  " }
    40 = @{ D = $null; E = "  -0.81%  " }
    41 = @{ D = $null; E = "  +1.44%  " }
    42 = @{ D = $null; E = "  -6.24%  " }
    43 = @{ D = $null; E = "  -2.68%  " }
    44 = @{ D = "96.50"; E = "  +0.91%  " }
    45 = @{ D = "2.86"; E = "  -13.22%  " }
    46 = @{ D = "16.06"; E = "  -6.94%  " }
    47 = @{ D = "1.352.00"; E = "  +6.14%  " }
    48 = @{ D = "2.43"; E = "  -0.60%  " }
    49 = @{ D = $null; E = "  +6.10%  " }
    50 = @{ D = $null; E = "  +1.16%  " }
    51 = @{ D = "2.265.31"; E = "  +1.47%  " }
}

foreach ($r in $updates.Keys) {
    $row = $updates[$r]
    if ($row.D -ne $null) {
        $ws.Range("D$r").NumberFormat = "@"
        $ws.Range("D$r").Value = $row.D
    }
    if ($row.E -ne $null) {
        $ws.Range("E$r").Value = $row.E
    }
}
